# Slide 5 ("Insert Data In Database"): split the explanatory sentence
# after "insert() " into multiple runs so that the literal tokens
# true/false, true and false are rendered in bold, matching the target
# OOXML (the sentence text itself is unchanged).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shape = $s.Shapes.Item(1)
$tf = $shape.TextFrame
$tr = $tf.TextRange

# Locate paragraph 2 ("insert() metodu true/false qaytarır. ...")
$para = $tr.Paragraphs(2, 1)
$paraText = $para.Text
$marker = "metodu"
$relIdx = $paraText.IndexOf($marker)
$base = $para.Start + $relIdx

# Sub-ranges (0-based offsets from $base) that must become bold.
# Everything else in the sentence keeps its existing (non-bold) run.
$boldSpans = @(
    @{ Offset = 7;  Length = 11 },  # "true/false "
    @{ Offset = 56; Length = 4 },   # "true"
    @{ Offset = 71; Length = 5 }    # "false"
)

foreach ($span in $boldSpans) {
    $sub = $tr.Characters($base + $span.Offset, $span.Length)
    $sub.Font.Bold = $true
}
